# Update the "想去人数" (want-to-go count) figures that changed between
# the two most recent data pulls, on both the "展览" sheet and the
# "全部类型" sheet (which mirrors the same rows).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 522
    $ws.Range("F6").Value = 86
    $ws.Range("F7").Value = 736
}
